# "added exact rank to each category and improved tie method"
#
# - Sheet1 (batting stats): a handful of "N" (at-bats?) values were
#   corrected/recomputed, and three tied rows (12, 13, 15) were given a new
#   underlined style to flag the tie-break.
# - The active sheet/selection moves from Sheet2 back to Sheet1.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# ---- Sheet1 data corrections -------------------------------------------
$ws1.Range("B4").Value  = 682
$ws1.Range("B5").Value  = 701
$ws1.Range("B10").Value = 682

# Tied rows get new values and a distinguishing style (underlined black
# Calibri 11, wrap text) to mark how the tie was broken.
foreach ($addr in @("B12", "B13", "B15")) {
    $cell = $ws1.Range($addr)
    $cell.Value = 660
    $cell.Font.Underline = $true
    $cell.Font.Color = 0
    $cell.WrapText = $true
}

# ---- Active sheet / selection moves back to Sheet1 ----------------------
[void]$ws1.Activate()
[void]$ws1.Range("K15").Select()
[void]$ws2.Range("J7").Select()
[void]$ws1.Activate()

Write-Host "done"
